$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.523.76"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "1.813.34"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("E16").Value = "  -5.38%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E18").Value = "  -3.63%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "26.588.29"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("E29").Value = "  -3.66%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("E51").Value = "  -1.53%  "

# Cells whose new text value parses as a plain number: force text storage
# without leaving a residual custom style on the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4566"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07139"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8814"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.301"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.377"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008609"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.019"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.983"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.080"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.863"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08695"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.043"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7327"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.674"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01963"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05130"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.894"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.995"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5009"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1556"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4608"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.972"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06003"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.54"
$ws.Range("D51").Style = "Normal"
